# --------------------------------------------------------------------------
# "break out stock.yaml completed"
#
# 1) The "day" sheet had five bsecode values (column D, rows 240-244) stored
#    as text; they should be plain numbers instead.
# 2) The "week" sheet gets fourteen newly scraped rows appended (103-116).
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) day sheet: convert bsecode text to numbers --------------------------
$wsDay = $wb.Worksheets.Item("day")

$wsDay.Cells.Item(240,4).Value = 500387
$wsDay.Cells.Item(241,4).Value = 532977
$wsDay.Cells.Item(242,4).Value = 533150
$wsDay.Cells.Item(243,4).Value = 531344
$wsDay.Cells.Item(244,4).Value = 500253

# --- 2) week sheet: append new rows -----------------------------------------
$wsWeek = $wb.Worksheets.Item("week")

function Set-WeekRow {
    param($row, $sr, $nsecode, $name, $bsecode, $perChg, $close, $volume, $timeframe, $dateTime)

    $wsWeek.Cells.Item($row, 1).Value = $sr
    $wsWeek.Cells.Item($row, 2).Value = $nsecode
    $wsWeek.Cells.Item($row, 3).Value = $name
    # bsecode is kept as text (like the rest of the sheet's existing rows),
    # so force it past numeric auto-detection with a leading apostrophe.
    $wsWeek.Cells.Item($row, 4).Value = "'" + $bsecode
    $wsWeek.Cells.Item($row, 5).Value = $perChg
    $wsWeek.Cells.Item($row, 6).Value = $close
    $wsWeek.Cells.Item($row, 7).Value = $volume
    $wsWeek.Cells.Item($row, 8).Value = $timeframe
    $wsWeek.Cells.Item($row, 9).Value = $dateTime
}

Set-WeekRow 103 1  "BOSCHLTD"   "Bosch Limited"                                   "500530" 0.15  34978.45   32525     "week" "31/07/2024 18:34:20"
Set-WeekRow 104 2  "POLYCAB"    "Polycab India Ltd"                               "542652" 2.64  6858.2     806650    "week" "31/07/2024 18:34:20"
Set-WeekRow 105 3  "JKCEMENT"   "Jk Cement Limited"                               "532644" -1.99 4416.95    126772    "week" "31/07/2024 18:34:20"
Set-WeekRow 106 4  "CUMMINSIND" "Cummins India Limited"                          "500480" 1.18  3852.35    386818    "week" "31/07/2024 18:34:20"
Set-WeekRow 107 5  "SHRIRAMFIN" "Shriram Finance Ltd"                            "511218" 0.68  2931.95    931158    "week" "31/07/2024 18:34:20"
Set-WeekRow 108 6  "BHARTIARTL" "Bharti Airtel Limited"                          "532454" 1.44  1491.55    5634305   "week" "31/07/2024 18:34:20"
Set-WeekRow 109 7  "IPCALAB"    "Ipca Laboratories Limited"                       "524494" 1.3   1307.95    335687    "week" "31/07/2024 18:34:20"
Set-WeekRow 110 8  "RAMCOCEM"   "The Ramco Cements Limited"                       "500260" 0.29  826.05     668598    "week" "31/07/2024 18:34:20"
Set-WeekRow 111 9  "AUBANK"     "AU Small Finance Bank"                           "540611" -0.83 646.05     3723813   "week" "31/07/2024 18:34:20"
Set-WeekRow 112 10 "INDHOTEL"   "The Indian Hotels Company Limited"               "500850" -0.64 642.05     1874629   "week" "31/07/2024 18:34:20"
Set-WeekRow 113 11 "EXIDEIND"   "Exide Industries Limited"                        "500086" -1.58 523.3      6871718   "week" "31/07/2024 18:34:20"
Set-WeekRow 114 12 "M&MFIN"     "Mahindra & Mahindra Financial Services Limited"  "532720" -0.25 303.25     1090343   "week" "31/07/2024 18:34:20"
Set-WeekRow 115 13 "RBLBANK"    "Rbl Bank Limited"                                "540065" -0.82 235.2      8441974   "week" "31/07/2024 18:34:20"
Set-WeekRow 116 14 "IDEA"       "Idea Cellular Limited"                           "532822" 0.37  16.27      346265142 "week" "31/07/2024 18:34:20"

# The apostrophe-prefix trick that forces text entry also stamps a
# "quote prefix" number format onto the cell; strip that back off so the new
# bsecode cells end up looking like ordinary (unstyled) text cells, matching
# the rest of the sheet.
$wsWeek.Range("D103:D116").ClearFormats()
